# edit.ps1 - applies the "Version beta poste 1" changes to Mon univers.docx
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: paragraph about DS6 Ruby Wild -- tweak a clause in the middle
# ---------------------------------------------------------------------------
$searchRange = $d.Content
$searchRange.Find.ClearFormatting()
$found1 = $searchRange.Find.Execute('DS6 Ruby Wild')
if (-not $found1) {
    throw "Could not locate anchor text 'DS6 Ruby Wild'"
}
$para1 = $searchRange.Paragraphs(1)
$para1.Range.Text = 'En 2011, j’ai participé à la conception, l’étude et à la réalisation du projet DS6 Ruby Wild destiné au marché chinois. Le rôle de concepteur dans un bureau de méthode permet de suivre la validation de chaque jalon clés d’un projet véhicule. Une expérience qui demande un savoir-faire technique, organisationnel et managérial dans un environnement pluridisciplinaire.'

# ---------------------------------------------------------------------------
# Edit 2: following paragraph about Megane/Clio/Alpine -- rewrite the sentence
# ---------------------------------------------------------------------------
$searchRange2 = $d.Content
$searchRange2.Find.ClearFormatting()
$found2 = $searchRange2.Find.Execute('Clio 4')
if (-not $found2) {
    throw "Could not locate anchor text 'Clio 4'"
}
$para2 = $searchRange2.Paragraphs(1)
$para2.Range.Text = 'Cette première expérience m’a permis de réaliser d’autre projet de grande envergure tel que l’étude de la Mégane 3, Clio 4 ainsi que le projet Renault Alpine, sortie en 2017.'

# ---------------------------------------------------------------------------
# Edit 3: replace the two "management" paragraphs under the Management
# heading with the much longer "Version beta poste 1" content (management,
# innovation quote, design-thinking paragraphs, etc.)
# ---------------------------------------------------------------------------
$searchRange3a = $d.Content
$searchRange3a.Find.ClearFormatting()
$found3a = $searchRange3a.Find.Execute('Le management a plusieurs aspects')
if (-not $found3a) {
    throw "Could not locate anchor text 'Le management a plusieurs aspects'"
}
$paraStart = $searchRange3a.Paragraphs(1)

$searchRange3b = $d.Content
$searchRange3b.Find.ClearFormatting()
$found3b = $searchRange3b.Find.Execute('Quelques soit le type de mangement')
if (-not $found3b) {
    throw "Could not locate anchor text 'Quelques soit le type de mangement'"
}
$paraEnd = $searchRange3b.Paragraphs(1)

$target = $d.Range($paraStart.Range.Start, $paraEnd.Range.End)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>En 2015, j’ai été emmené à manager une équipe de « design to cost » conception dont l’objectif était de mener des études en avance de phase sur 25 projets automobile. Piloter une équipe d’expertise métier dans un but commun a été pour moi, une expérience particulièrement enrichissante. De part, la gestion des objectifs, j’ai compris l’importance de mettre en premier plan la gestion des membres de mon équipes.</w:t></w:r></w:p><w:p><w:r><w:t>Il ne faut pas oublier que la valeur du produit est avant tout réalisé par les équipes.</w:t></w:r></w:p><w:p><w:r><w:t>Le management dit</w:t></w:r><w:r><w:t xml:space="preserve"> innovant tel que le management 3.0 permet d''augmenter la performance des équipes à travers la création de challenge permanent</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Innovation</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>« L''innovation systématique requiert la volonté de considérer le changement comme une opportunité. »</w:t></w:r></w:p><w:p><w:r><w:t>Peter Drucker - Artiste, écrivain, Enseignant (1909 - 2005)</w:t></w:r></w:p><w:p><w:r><w:t>À travers mes expériences dans les différents secteurs d’activité, j’ai relevé un seul moins point commun : la nécessité d’innover pour subsister.</w:t></w:r></w:p><w:p><w:r><w:t>L’innovation est devenue la clé de voute pour toutes les entreprises que souhaite rester sur un marché de plus en plus exigent.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">En tant que coach agile, j’ai été emmené à construire des ateliers de brainstormings dans le but d’augmenter la créativité et la productivité des équipes. Depuis quelques années, la démarche de design </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>thinking</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> est de plus en plus utilisé pour favoriser l’innovation.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml)
